# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 09:05"

# Row 121/122: Georgia and Crucero swap places (Georgia now listed before Crucero)
# and each gets refreshed case data.
$ws.Range("A121").Value = "Georgia"
$ws.Range("B121").Value = 713
$ws.Range("C121").Value = 6
$ws.Range("D121").Value = 475
$ws.Range("E121").Value = 226
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 12

$ws.Range("A122").Value = "Crucero"
$ws.Range("B122").Value = 712
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 651
$ws.Range("E122").Value = 48
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 13

# Row 66: Armenia
$ws.Range("B66").Value = 5271
$ws.Range("C66").Value = 230
$ws.Range("D66").Value = 2419
$ws.Range("E66").Value = 2785
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 67

# Row 106: Letonia
$ws.Range("B106").Value = 1016
$ws.Range("C106").Value = 4
$ws.Range("D106").Value = 694
$ws.Range("E106").Value = 301
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 21

# Row 132: Taiwan
$ws.Range("D132").Value = 402
$ws.Range("E132").Value = 31
